$p = $ppt.ActivePresentation

# --- Modify slide 8 ("Card") content: merge the "color, id" bullet with the
# "Construct a null..." bullet into a single "color, number, function, id" bullet.
$s8 = $p.Slides.Item(8)
$contentShape = $s8.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

# Delete the second paragraph ("Construct a null to represent the situation of "no card"")
$para2 = $tr.Paragraphs(2, 1)
$para2.Delete()

# Replace the remaining first paragraph's text. Assign an intermediate
# placeholder first so the engine doesn't try to diff/split the old and new
# text into multiple runs, keeping a single run like the authored slide.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "-"
$para1b = $tr.Paragraphs(1, 1)
$para1b.Text = "Store and get the color, number, function, id of the card"

# --- Remove the last two slides (NumberCard, FunctionCard) ---
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()
